# Update the data model: rename the "has*"-prefixed technical header labels
# in row 1 to their human-readable equivalents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "File Name"
$ws.Range("C1").Value = "Time Stamp"
$ws.Range("D1").Value = "Copyright"
$ws.Range("E1").Value = "License List"
$ws.Range("F1").Value = "Date"

# Update the view state: zoom level and active selection.
$excel.ActiveWindow.Zoom = 207
$ws.Range("G7").Select() | Out-Null
